$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# "Floppy Maker (unity)": D7 currently shows the "kb" label; change it to
# a new, distinct label "b" (adds a new shared string, since no existing
# string in the sheet currently has that value).
$ws.Range("D7").Value = "b"

# Move the sheet's selected/active cell from G9 to C7.
$ws.Range("C7").Select()

# "Floppy Interrupt": minimize the workbook window.
$excel.ActiveWindow.WindowState = -4140
